# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block
# (and the blank paragraph that separates it from the bibliography),
# leaving the bibliography's last entry (SAMUELSON ...) followed directly
# by the existing blank paragraph + page-break paragraph that close the
# document.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph.
$jupiterPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $jupiterPara = $p
        break
    }
}

# Locate the "(c) 2020 ... Creative Commons Attribution" paragraph.
$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Powered by Jekyll and Github pages*") {
        $copyrightPara = $p
        break
    }
}

# The blank paragraph immediately preceding the "Ver no Jupiter ..." one
# (sits right after the SAMUELSON bibliography entry) is removed too.
$blankPara = $jupiterPara.Previous(1)

$startRange = $blankPara.Range.Start
$endRange = $copyrightPara.Range.End

$toDelete = $d.Range($startRange, $endRange)
$toDelete.Delete()
